# Configuration_Audit.xlsx - add Maria DB / MSSQL rows, fix MSSQL naming
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing Database rows (B6:B8) ---
# Row 6: "Microsoft SQL Server 2019 Benchmark v1.5.0" -> "Microsoft SQL Server 2019 " (keep trailing space)
$ws.Range("B6").Value = "Microsoft SQL Server 2019 "

# Row 7: "Microsoft SQL Server 2017 Benchmark v1.3.0" -> "Microsoft SQL Server 2017 " (keep trailing space)
$ws.Range("B7").Value = "Microsoft SQL Server 2017 "

# Row 8: "Microsoft SQL Server 2016 Benchmark v1.4.0" -> "Microsoft SQL Server 2022"
$ws.Range("B8").Value = "Microsoft SQL Server 2022"

# --- Append new Database rows 9-11, matching formatting of row 7/8 (style index "1", row height 15) ---
$ws.Range("B7").Copy()
$ws.Range("B9:B11").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Rows.Item(9).RowHeight = $ws.Rows.Item(7).RowHeight
$ws.Rows.Item(10).RowHeight = $ws.Rows.Item(7).RowHeight
$ws.Rows.Item(11).RowHeight = $ws.Rows.Item(7).RowHeight

# Row 9: Microsoft SQL Server 2016 (trailing space)
$ws.Range("A9").Value = "Database"
$ws.Range("B9").Value = "Microsoft SQL Server 2016 "
$ws.Range("C9").Value = "CIS"

# Row 10: Maria Server v10_6
$ws.Range("A10").Value = "Database"
$ws.Range("B10").Value = "Maria Server v10_6"
$ws.Range("C10").Value = "CIS"

# Row 11: Maria Server v10_11
$ws.Range("A11").Value = "Database"
$ws.Range("B11").Value = "Maria Server v10_11"
$ws.Range("C11").Value = "CIS"

# --- Update the hidden _FilterDatabase defined name range (A1:A66 -> A1:A67) ---
$n = $wb.Names.Item("_xlnm._FilterDatabase")
$n.RefersTo = "=Sheet1!`$A`$1:`$A`$67"
